$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text so numeric-looking values are not
#     auto-converted to numbers (matches source data stored as text) ---
$priceCells = @("D2","D3","D5","D6","D8","D9","D13","D14","D15","D16","D18","D19","D20","D22","D24","D26","D27","D28","D30","D32","D33","D36","D37","D38","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.230.68"
$ws.Range("D3").Value = "2.432.20"
$ws.Range("D5").Value = "565.59"
$ws.Range("D6").Value = "141.20"
$ws.Range("D8").Value = "0.590"
$ws.Range("D9").Value = "2.431.98"
$ws.Range("D13").Value = "0.354"
$ws.Range("D14").Value = "26.24"
$ws.Range("D15").Value = "2.869.28"
$ws.Range("D16").Value = "63.057.12"
$ws.Range("D18").Value = "2.428.91"
$ws.Range("D19").Value = "11.25"
$ws.Range("D20").Value = "340.90"
$ws.Range("D22").Value = "6.82"
$ws.Range("D24").Value = "65.27"
$ws.Range("D26").Value = "1.00"
$ws.Range("D27").Value = "1.53"
$ws.Range("D28").Value = "8.17"
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("D32").Value = "6.63"
$ws.Range("D33").Value = "174.01"
$ws.Range("D36").Value = "18.71"
$ws.Range("D37").Value = "369.83"
$ws.Range("D38").Value = "4.48"
$ws.Range("D41").Value = "1.69"
$ws.Range("D42").Value = "39.94"
$ws.Range("D43").Value = "148.50"
$ws.Range("D44").Value = "3.69"
$ws.Range("D45").Value = "20.76"
$ws.Range("D46").Value = "0.594"
$ws.Range("D47").Value = "0.0958"
$ws.Range("D48").Value = "0.0522"
$ws.Range("D49").Value = "0.0224"
$ws.Range("D50").Value = "17.87"
$ws.Range("D51").Value = "1.73"

# --- Volume(1h) column (E): percentage text values ---
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("E3").Value = "  +5.51%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("E6").Value = "  +9.85%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("E10").Value = "  +4.22%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +7.04%  "
$ws.Range("E14").Value = "  +13.37%  "
$ws.Range("E15").Value = "  +5.67%  "
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("E17").Value = "  +8.78%  "
$ws.Range("E18").Value = "  +5.92%  "
$ws.Range("E19").Value = "  +8.02%  "
$ws.Range("E20").Value = "  +9.88%  "
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("E22").Value = "  +4.77%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +13.60%  "
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("E29").Value = "  +11.96%  "
$ws.Range("E30").Value = "  +11.93%  "
$ws.Range("E31").Value = "  +7.26%  "
$ws.Range("E32").Value = "  +14.32%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  +11.68%  "
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("E37").Value = "  +17.87%  "
$ws.Range("E38").Value = "  +12.32%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +12.74%  "
$ws.Range("E42").Value = "  +6.58%  "
$ws.Range("E43").Value = "  +9.08%  "
$ws.Range("E44").Value = "  +8.16%  "
$ws.Range("E45").Value = "  +11.77%  "
$ws.Range("E46").Value = "  +4.71%  "
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  +6.85%  "
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("E50").Value = "  +7.50%  "
$ws.Range("E51").Value = "  +16.15%  "
